# Actualización automática de noticias - 2026-01-14
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (row 2) to make room for the
# two newest news items, pushing all existing data rows down by two rows.
$ws.Range("A2:F3").EntireRow.Insert()

# Row 2: new item (2026-01-14) - Infobae
# (dates are plain text in this sheet; a leading apostrophe forces text
# entry so Excel doesn't reinterpret "2026-01-14" as a date serial)
$ws.Cells.Item(2, 1).Value = "'2026-01-14"
$ws.Cells.Item(2, 2).Value = "Menor que fue secuestrada en el Catatumbo denunció que las disidencias la obligaron a enviar mensaje de terror a otros niños"
$ws.Cells.Item(2, 3).Value = "Infobae"
$ws.Cells.Item(2, 4).Value = "Sin identificar"
$ws.Cells.Item(2, 5).Value = "https://www.infobae.com/colombia/2026/01/14/menor-que-fue-secuestrada-en-el-catatumbo-denuncio-que-las-disidencias-la-obligaron-a-enviar-mensaje-de-terror-a-otros-ninos/"
$ws.Cells.Item(2, 6).Value = "PorJhon Bernal"

# Row 3: new item (2026-01-14) - Infobae
$ws.Cells.Item(3, 1).Value = "'2026-01-14"
$ws.Cells.Item(3, 2).Value = "Menor de edad colombiano fue baleado frente a su escuela en Chicago, la familia exige justicia"
$ws.Cells.Item(3, 3).Value = "Infobae"
$ws.Cells.Item(3, 4).Value = "Colombia"
$ws.Cells.Item(3, 5).Value = "https://www.infobae.com/colombia/2026/01/14/menor-de-edad-colombiano-fue-baleado-frente-a-su-escuela-en-chicago-la-familia-exige-justicia/"
$ws.Cells.Item(3, 6).Value = "PorJimmy Nomesqui Rivera"

# The insert above copied the header row's bold/bordered style into the two
# new rows; strip that back to the plain (unstyled) look the other data
# rows use.
$ws.Range("A2:F3").ClearFormats()
